# Apply cryptos list update (GitHub Actions refresh) to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" strings look numeric (e.g. "290.03"). Assigning them
# directly would let Excel coerce the text into a floating point
# number, so for those cells we briefly force Text format, write the
# literal string, then restore each cells original (unstyled) look.

$ws.Range("D2").Value = "40.046.18"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "2.224.17"
$ws.Range("E3").Value = "  -0.42%  "

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.07%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "290.03"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -1.19%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.29"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +3.77%  "

$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("E8").Value = "  -0.08%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.473"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +0.81%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.62"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +1.04%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -2.12%  "

$ws.Range("E12").Value = "  +2.96%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.52"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +2.48%  "

$ws.Range("D14").Value = "2.563.83"
$ws.Range("E14").Value = "  -0.43%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.00"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").Value = "2.214.16"
$ws.Range("E16").Value = "  -0.57%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.731"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "39.987.63"
$ws.Range("E18").Value = "  +0.43%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.59"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +9.79%  "

$ws.Range("D20").Value = "0.0₃0884"
$ws.Range("E20").Value = "  -0.79%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.81"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +0.39%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.70"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.45%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.28"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +1.77%  "

$ws.Range("E24").Value = "  -0.03%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +1.70%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.83"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.72%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.61"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -1.49%  "

$ws.Range("E28").Value = "  -0.14%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.24"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +0.15%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "155.95"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +0.49%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.92"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -3.11%  "

$ws.Range("E32").Value = "  -0.10%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.95"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +2.25%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0719"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +1.77%  "

$ws.Range("E35").Value = "  +0.63%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.87"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +7.53%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.112"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -0.04%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.87"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -3.97%  "

$ws.Range("E39").Value = "  +0.50%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.71"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +2.92%  "

$ws.Range("D41").Value = "2.109.76"
$ws.Range("E41").Value = "  +8.15%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.85"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +3.14%  "

$ws.Range("E43").Value = "  -1.64%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0268"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.95"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +6.16%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.48"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +7.72%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.67"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +2.66%  "

$ws.Range("D48").Value = "2.433.25"
$ws.Range("E48").Value = "  -0.50%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.45"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "88.91"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "69.22"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -2.32%  "

Write-Output "cryptos list updated"